$wb = $excel.ActiveWorkbook

### Sheet: 展览 ###
$ws = $wb.Worksheets.Item('展览')

# Update 'want to go' counts (column F)
$ws.Range('F3').Value = 184
$ws.Range('F4').Value = 181
$ws.Range('F5').Value = 5096
$ws.Range('F8').Value = 17
$ws.Range('F13').Value = 1420
$ws.Range('F14').Value = 3767
$ws.Range('F15').Value = 419
$ws.Range('F16').Value = 147
$ws.Range('F17').Value = 133
$ws.Range('F19').Value = 2864
$ws.Range('F20').Value = 139
$ws.Range('F21').Value = 178
$ws.Range('F22').Value = 94
$ws.Range('F25').Value = 74
$ws.Range('F26').Value = 18
$ws.Range('F29').Value = 281

# Update venue text for row 4
$ws.Range('D4').Value = '滨江大会堂 滨江宾馆'

# Insert a new row for the added '信丰·端午节UPUP动漫展' entry,
# shifting the existing rows (formerly starting at row 30) down by one
$ws.Rows.Item(30).Insert()
# Restore the column-A cell format (border/bold/center) that Insert()
# does not carry over identically, by re-using the row above's format
$ws.Range('A29').Copy()
$ws.Range('A30').PasteSpecial(-4122)

# Populate the newly-inserted row 30
$ws.Range('A30').Value = 29
$ws.Range('B30').Value = '2024-06-09'
$ws.Range('C30').Value = '信丰·端午节UPUP动漫展'
$ws.Range('D30').Value = '迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子运动馆'
$ws.Range('E30').Value = '2024.06.09 10:00-06.09 17:00'
$ws.Range('F30').Value = 0
$ws.Range('G30').Value = 48
$ws.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=84078'
$ws.Range('I30').Value = '//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg'

### Sheet: 全部类型 ###
$ws = $wb.Worksheets.Item('全部类型')

# Update 'want to go' counts (column F)
$ws.Range('F3').Value = 184
$ws.Range('F4').Value = 181
$ws.Range('F6').Value = 5096
$ws.Range('F9').Value = 17
$ws.Range('F14').Value = 1420
$ws.Range('F15').Value = 3767
$ws.Range('F16').Value = 419
$ws.Range('F17').Value = 147
$ws.Range('F18').Value = 133
$ws.Range('F20').Value = 2864
$ws.Range('F21').Value = 139
$ws.Range('F22').Value = 178
$ws.Range('F23').Value = 94
$ws.Range('F26').Value = 74
$ws.Range('F27').Value = 18
$ws.Range('F30').Value = 281

# Update venue text for row 4
$ws.Range('D4').Value = '滨江大会堂 滨江宾馆'

# Insert a new row for the added '信丰·端午节UPUP动漫展' entry,
# shifting the existing rows (formerly starting at row 31) down by one
$ws.Rows.Item(31).Insert()
# Restore the column-A cell format (border/bold/center) that Insert()
# does not carry over identically, by re-using the row above's format
$ws.Range('A30').Copy()
$ws.Range('A31').PasteSpecial(-4122)

# Populate the newly-inserted row 31
$ws.Range('A31').Value = 30
$ws.Range('B31').Value = '2024-06-09'
$ws.Range('C31').Value = '信丰·端午节UPUP动漫展'
$ws.Range('D31').Value = '迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子运动馆'
$ws.Range('E31').Value = '2024.06.09 10:00-06.09 17:00'
$ws.Range('F31').Value = 0
$ws.Range('G31').Value = 48
$ws.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=84078'
$ws.Range('I31').Value = '//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg'

